# Graduation students template:
#   - add a "Căn cước công dân" (citizen identification) column
#   - update a couple of data values
#   - refresh the table's visual formatting (font, borders, wrap)
#     to match how the sheet looked after being re-saved/tidied in Excel

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Content edits first (so they ride along when formatting is applied) --

# New column header + data: citizen identification number
$ws.Range("H1").Value = "Căn cước công dân"

# Mã sinh viên (student code) changed
$ws.Range("A2").Value = 622222

# Điểm trung bình (GPA) becomes a real number instead of text "3.6"
$ws.Range("D2").Value = 3.6

# Xếp Loại (classification) cell picks up a stray newline-only value
$ws.Range("E2").Value = "`n"

# ---- Visual formatting: font + borders across the whole table -------------
$tbl = $ws.Range("A1:H2")
$tbl.Font.Name = "Helvetica"
$tbl.Font.Size = 8

$tbl.Borders.Color = 12566463
$tbl.Borders.LineStyle = 1

$ws.Range("H1").WrapText = $true
$ws.Range("E2").WrapText = $true

# H2 last: keep the leading zero by forcing text, then wrap
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0123456789"
$ws.Range("H2").WrapText = $true

# ---- Column widths / row heights ------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 11.43
$ws.Columns.Item(7).ColumnWidth = 15
$ws.Columns.Item(8).ColumnWidth = 14

$ws.Rows.Item(1).RowHeight = 23.25
$ws.Rows.Item(2).RowHeight = 23.25

# ---- Selection --------------------------------------------------------
$ws.Range("F8").Select()
